$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (values recomputed by the
# quarterly-indexing bug fix).
$data = @{
    2  = @(0.433935198812251,  0.4883748341653806, 0.2913197292499859, 0.5397404276594314, 0.3340733006221409, 13)
    3  = @(0.3043483804216044, 0.3187202035556604, 0.1391823142959026, 0.3730714600393638, 0.2253588111853846, 12)
    4  = @(0.2239330230754586, 0.2420726978122317, 0.08325831886165148, 0.2885451764657512, 0.190849553422929, 11)
    5  = @(0.3439814776516536, 0.3439814776516536, 0.1481954533256283, 0.3849616257831789, 0.1833199959169479, 9)
    6  = @(0.3607780907168635, 0.3638782949054072, 0.166400635697807,  0.40792234027791,   0.2006650192526551, 10)
    7  = @(0.3167374129797205, 0.3281249448458587, 0.1360088531859135, 0.3687937813818361, 0.2003672813995112, 9)
    8  = @(0.3332309007927744, 0.3363297904028241, 0.1472445396808239, 0.3837245622589515, 0.2034044709093121, 8)
    9  = @(0.2929253822463562, 0.295209788697302,  0.1182577732422086, 0.3438862795201469, 0.1973397892307624, 6)
    10 = @(0.3047033996258702, 0.3047033996258702, 0.1225151295288574, 0.35002161294534,   0.1886932996753026, 6)
    11 = @(0.2886753471776797, 0.2886753471776797, 0.1109408868566417, 0.3330778990816439, 0.1857667582900928, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
